# The dataset was renamed: the "congenital" variable/category entry found in
# column A of each "variables_*" worksheet is renamed to "misc_long_term".
# Not every worksheet contains this value, so each sheet is searched and only
# updated when a match is found.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $found = $ws.UsedRange.Find("congenital")
    if ($found -ne $null) {
        $found.Value = "misc_long_term"
    }
}
